$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title (it is being moved/replaced elsewhere).
# ------------------------------------------------------------------
foreach ($metaPara in $d.Paragraphs) {
    if ($metaPara.Range.Text.StartsWith("Meta description")) {
        $metaPara.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph, "Play A While On The Nile Slot for
#    Free - Exciting Bonus Features", right after the last bullet
#    point ("May require some investment to land noteworthy wins")
#    and before the final italic paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$bulletPara = $d.Paragraphs.Item($count - 1)
$bulletPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newPara.Style = "Normal"

$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertPoint.InsertAfter("Play A While On The Nile Slot for Free - Exciting Bonus Features")

$newPara2 = $d.Paragraphs.Item($count)
$newParaText = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
$newParaText.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Replace the old "Create a Feature Image Prompt: ..." text (the
#    final, italic paragraph) with the new meta-description copy.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a Feature Image Prompt: Design a cartoon-style feature image for " + [char]34 + "A While On The Nile" + [char]34 + " online slot game that showcases a happy Maya warrior wearing glasses. The warrior should be depicted in an Egyptian-themed outfit, possibly holding a tablet with hieroglyphics or standing by the Nile river. The image should have bright colors and should be eye-catching to potential players. Make sure to include the game's title in the image prominently.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experience the ancient Egypt theme with A While On The Nile and enjoy exciting bonus features - play for free and potentially win big.",
    2
) | Out-Null
